$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date column C for all existing data rows
#    (rows 2-504) from 45190 (2023-09-21) to 45192 (2023-09-23).
$ws.Range("C2:C504").Value = 45192

# 2. Row 504 gains an explicit row height (15pt, custom height).
$ws.Rows.Item(504).RowHeight = 15

# 3. Append a new record as row 505.
$ws.Cells.Item(505, 1).Value = "A 44712-2023"

$ws.Cells.Item(505, 2).Value = 45190
$ws.Cells.Item(505, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(505, 3).Value = 45192
$ws.Cells.Item(505, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(505, 4).Value = "DALARNAS LÄN"
$ws.Cells.Item(505, 5).Value = "LEKSAND"
$ws.Cells.Item(505, 6).Value = "Bergvik skog väst AB"

$ws.Cells.Item(505, 7).Value = 1.2
$ws.Cells.Item(505, 8).Value = 0
$ws.Cells.Item(505, 9).Value = 0
$ws.Cells.Item(505, 10).Value = 0
$ws.Cells.Item(505, 11).Value = 0
$ws.Cells.Item(505, 12).Value = 0
$ws.Cells.Item(505, 13).Value = 0
$ws.Cells.Item(505, 14).Value = 0
$ws.Cells.Item(505, 15).Value = 0
$ws.Cells.Item(505, 16).Value = 0
$ws.Cells.Item(505, 17).Value = 0

# Column R ("Artnamn") is left as an empty, wrap-text styled cell like the
# other rows that have no species findings.
$ws.Cells.Item(505, 18).Value = "x"
$ws.Cells.Item(505, 18).WrapText = $true
$ws.Cells.Item(505, 18).Value = ""
